$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 5734.1665
$ws.Range("I5").Value = 450
$ws.Range("J5").Value = 6791
$ws.Range("K5").Value = 450
$ws.Range("L5").Value = 6791
$ws.Range("M5").Value = -335
$ws.Range("N5").Value = -7021
$ws.Range("H15").Value = 307238.53
$ws.Range("I15").Value = 307238.53
$ws.Range("K15").Value = 921715.5900000001
$ws.Range("M15").Value = -921546.5900000001
$ws.Range("H31").Value = 915.4
$ws.Range("I31").Value = 1106.75
$ws.Range("J31").Value = 150
$ws.Range("K31").Value = 3320.25
$ws.Range("L31").Value = 450
$ws.Range("M31").Value = -3090.25
$ws.Range("N31").Value = -910
$ws.Range("H33").Value = 842.5789
$ws.Range("I33").Value = 794.9286
$ws.Range("J33").Value = 976
$ws.Range("K33").Value = 794.9286
$ws.Range("L33").Value = 976
$ws.Range("M33").Value = -565.9286
$ws.Range("N33").Value = -1434
$ws.Range("H76").Value = 5992
$ws.Range("I76").Value = 5365.1333
$ws.Range("J76").Value = 7335.2856
$ws.Range("K76").Value = 5365.1333
$ws.Range("L76").Value = 7335.2856
$ws.Range("M76").Value = -5050.1333
$ws.Range("N76").Value = -7965.2856
$ws.Range("H79").Value = 5992
$ws.Range("I79").Value = 5365.1333
$ws.Range("J79").Value = 7335.2856
$ws.Range("K79").Value = 5365.1333
$ws.Range("L79").Value = 7335.2856
$ws.Range("M79").Value = -4273.1333
$ws.Range("N79").Value = -9519.285599999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 3213.5117
$ws.Range("I32").Value = 1265.6
$ws.Range("K32").Value = 1265.6
$ws.Range("M32").Value = -978.5999999999999
$ws.Range("H61").Value = 6388.75
$ws.Range("I61").Value = 8912.4
$ws.Range("K61").Value = 8912.4
$ws.Range("M61").Value = -8700.4
$ws.Range("H74").Value = 6422
$ws.Range("I74").Value = 8005.5
$ws.Range("J74").Value = 5366.3335
$ws.Range("K74").Value = 8005.5
$ws.Range("L74").Value = 5366.3335
$ws.Range("M74").Value = -7131.5
$ws.Range("N74").Value = -7114.3335
$ws.Range("H77").Value = 6422
$ws.Range("I77").Value = 8005.5
$ws.Range("J77").Value = 5366.3335
$ws.Range("K77").Value = 40027.5
$ws.Range("L77").Value = 26831.6675
$ws.Range("M77").Value = -35659.5
$ws.Range("N77").Value = -35567.6675
$ws.Range("H94").Value = 59250
$ws.Range("J94").Value = 59250
$ws.Range("L94").Value = 59250
$ws.Range("N94").Value = -61052
$ws.Range("H110").Value = 4895.2334
$ws.Range("J110").Value = 6180
$ws.Range("L110").Value = 6180
$ws.Range("N110").Value = -10270
$ws.Range("H132").Value = 70499.71000000001
$ws.Range("I132").Value = 4030.75
$ws.Range("K132").Value = 12092.25
$ws.Range("M132").Value = -9562.25
$ws.Range("H136").Value = 6388.75
$ws.Range("I136").Value = 8912.4
$ws.Range("K136").Value = 26737.2
$ws.Range("M136").Value = -24187.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1761.1666
$ws.Range("J20").Value = 2175.8462
$ws.Range("L20").Value = 2175.8462
$ws.Range("N20").Value = -2669.8462
$ws.Range("H22").Value = 2821.8965
$ws.Range("I22").Value = 1417.9445
$ws.Range("J22").Value = 5119.273
$ws.Range("K22").Value = 1417.9445
$ws.Range("L22").Value = 5119.273
$ws.Range("M22").Value = -1244.9445
$ws.Range("N22").Value = -5465.273
$ws.Range("H99").Value = 8646.375
$ws.Range("I99").Value = 4587.5
$ws.Range("J99").Value = 9999.333000000001
$ws.Range("K99").Value = 4587.5
$ws.Range("L99").Value = 9999.333000000001
$ws.Range("M99").Value = -3089.5
$ws.Range("N99").Value = -12995.333
$ws.Range("H134").Value = 8742.1
$ws.Range("I134").Value = 7782.8887
$ws.Range("J134").Value = 17375
$ws.Range("K134").Value = 23348.6661
$ws.Range("L134").Value = 52125
$ws.Range("M134").Value = -20813.6661
$ws.Range("N134").Value = -57195

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 6696.778
$ws.Range("I16").Value = 4454.4
$ws.Range("K16").Value = 4454.4
$ws.Range("M16").Value = -4167.4
$ws.Range("H31").Value = 20004
$ws.Range("I31").Value = 5006
$ws.Range("J31").Value = 50000
$ws.Range("K31").Value = 5006
$ws.Range("L31").Value = 50000
$ws.Range("M31").Value = -4711
$ws.Range("N31").Value = -50590
$ws.Range("H34").Value = 20004
$ws.Range("I34").Value = 5006
$ws.Range("J34").Value = 50000
$ws.Range("K34").Value = 5006
$ws.Range("L34").Value = 50000
$ws.Range("M34").Value = -4804
$ws.Range("N34").Value = -50404
$ws.Range("H39").Value = 6295.4165
$ws.Range("I39").Value = 7590.8335
$ws.Range("J39").Value = 5000
$ws.Range("K39").Value = 7590.8335
$ws.Range("L39").Value = 5000
$ws.Range("M39").Value = -7199.8335
$ws.Range("N39").Value = -5782
$ws.Range("H49").Value = 6295.4165
$ws.Range("I49").Value = 7590.8335
$ws.Range("J49").Value = 5000
$ws.Range("K49").Value = 7590.8335
$ws.Range("L49").Value = 5000
$ws.Range("M49").Value = -7408.8335
$ws.Range("N49").Value = -5364
$ws.Range("H105").Value = 11333
$ws.Range("I105").Value = 13333
$ws.Range("J105").Value = 9333
$ws.Range("K105").Value = 13333
$ws.Range("L105").Value = 9333
$ws.Range("M105").Value = -11586
$ws.Range("N105").Value = -12827
$ws.Range("H113").Value = 6696.778
$ws.Range("I113").Value = 4454.4
$ws.Range("K113").Value = 4454.4
$ws.Range("M113").Value = -2284.4
$ws.Range("H122").Value = 3890.25
$ws.Range("I122").Value = 3377.4285
$ws.Range("J122").Value = 4608.2
$ws.Range("K122").Value = 10132.2855
$ws.Range("L122").Value = 13824.6
$ws.Range("M122").Value = -7682.2855
$ws.Range("N122").Value = -18724.6
$ws.Range("H132").Value = 5387.9546
$ws.Range("I132").Value = 4622.4443
$ws.Range("K132").Value = 13867.3329
$ws.Range("M132").Value = -11337.3329

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H56").Value = 8486.076999999999
$ws.Range("I56").Value = 8486.076999999999
$ws.Range("K56").Value = 8486.076999999999
$ws.Range("M56").Value = -7956.076999999999
$ws.Range("H108").Value = 5030.1333
$ws.Range("I108").Value = 1439.2222
$ws.Range("K108").Value = 4317.6666
$ws.Range("M108").Value = -1437.6666

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H39").Value = 0
$ws.Range("J39").Value = 0
$ws.Range("L39").Value = 0
$ws.Range("H70").Value = 13759.857
$ws.Range("I70").Value = 6925.7
$ws.Range("K70").Value = 6925.7
$ws.Range("M70").Value = -6655.7
$ws.Range("H73").Value = 13759.857
$ws.Range("I73").Value = 6925.7
$ws.Range("K73").Value = 6925.7
$ws.Range("M73").Value = -5989.7
$ws.Range("H93").Value = 45250
$ws.Range("J93").Value = 45250
$ws.Range("L93").Value = 45250
$ws.Range("N93").Value = -48994
$ws.Range("H102").Value = 7171.9395
$ws.Range("I102").Value = 6458.0527
$ws.Range("J102").Value = 8140.7856
$ws.Range("K102").Value = 6458.0527
$ws.Range("L102").Value = 8140.7856
$ws.Range("M102").Value = -4836.0527
$ws.Range("N102").Value = -11384.7856
$ws.Range("H122").Value = 17195
$ws.Range("I122").Value = 19801.4
$ws.Range("K122").Value = 59404.2
$ws.Range("M122").Value = -56954.2
$ws.Range("H132").Value = 2504.5715
$ws.Range("I132").Value = 2504.5715
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 7513.7145
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -4983.7145
$ws.Range("H133").Value = 156331
$ws.Range("J133").Value = 156331
$ws.Range("L133").Value = 156331
$ws.Range("N133").Value = -166451
$ws.Range("H135").Value = 125999.336
$ws.Range("I135").Value = 109999.5
$ws.Range("J135").Value = 130570.71
$ws.Range("K135").Value = 109999.5
$ws.Range("L135").Value = 130570.71
$ws.Range("M135").Value = -104929.5
$ws.Range("N135").Value = -140710.71
$ws.Range("N132").ClearContents()
$ws.Range("N39").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6650.643
$ws.Range("I7").Value = 5020.125
$ws.Range("K7").Value = 5020.125
$ws.Range("M7").Value = -4908.125
$ws.Range("H16").Value = 3723.476
$ws.Range("I16").Value = 1625.5714
$ws.Range("J16").Value = 7919.2856
$ws.Range("K16").Value = 1625.5714
$ws.Range("L16").Value = 7919.2856
$ws.Range("M16").Value = -1455.5714
$ws.Range("N16").Value = -8259.285599999999
$ws.Range("H82").Value = 2571.2258
$ws.Range("I82").Value = 1408.5
$ws.Range("K82").Value = 1408.5
$ws.Range("M82").Value = -1047.5
$ws.Range("H85").Value = 2571.2258
$ws.Range("I85").Value = 1408.5
$ws.Range("K85").Value = 1408.5
$ws.Range("M85").Value = -160.5
$ws.Range("H126").Value = 6650.643
$ws.Range("I126").Value = 5020.125
$ws.Range("K126").Value = 15060.375
$ws.Range("M126").Value = -12590.375

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H46").Value = 102721.164
$ws.Range("J46").Value = 102721.164
$ws.Range("L46").Value = 102721.164
$ws.Range("N46").Value = -103183.164
$ws.Range("H122").Value = 2932.0667
$ws.Range("I122").Value = 2566.3333
$ws.Range("J122").Value = 3175.889
$ws.Range("K122").Value = 7698.999899999999
$ws.Range("L122").Value = 9527.667000000001
$ws.Range("M122").Value = -5248.999899999999
$ws.Range("N122").Value = -14427.667
$ws.Range("H132").Value = 8312.16
$ws.Range("I132").Value = 5524.25
$ws.Range("K132").Value = 16572.75
$ws.Range("M132").Value = -14042.75
$ws.Range("H134").Value = 102721.164
$ws.Range("J134").Value = 102721.164
$ws.Range("L134").Value = 308163.492
$ws.Range("N134").Value = -313233.492
